# Criação da tela de autorização para finaliza_tratamento.php
# Adds 22 new "GRANT ALL ON SEQUENCE ..." rows (185:206) for the new user
# "mariana.brider" to the "grants por usuario" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# Column A text (GRANT ... TO) for new rows 185-206, in row order so that
# any brand-new shared strings are appended to sharedStrings.xml in this
# exact sequence.
$aValues = @(
    "GRANT ALL ON SEQUENCE tratamento.sq_acesso_transac_tratamento TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_cnvo TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_equipe TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_equipe_usua TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_grupo_acesso TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_grupo_usua_acesso TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_grupo_usua_menu_sist_tratamento TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_grupo_usua_transac_acesso TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_hstr_obs_pnel_mapa_risco TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_hstr_obs_pnel_solic_trtmto TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_hstr_pnel_mapa_risco TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_hstr_pnel_solic_trtmto TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_local_trtmto TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_log_acesso TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_menu_sist_tratamento TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_pddo_trtmto TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_pnel_solic_trtmto TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_risco_pcnt TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_risco_rnado_pcnt TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_status_pcnt TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_status_trtmto TO",
    "GRANT ALL ON SEQUENCE tratamento.sq_usua_acesso TO"
)

$firstRow = 185
$lastRow = 206

# 1) Column A - fill every row first (matches original authoring order,
#    and is what drives the order new shared strings get created in).
for ($i = 0; $i -lt $aValues.Count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value = $aValues[$i]
}

# 2) Column B - the new user, all rows.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = '"mariana.brider"'
}
# Copy column B's styling (style index 1: Arial 9 / wrap / vertical center)
# from the row above, same as every other data row in the column.
$ws.Range("B184").Copy()
$ws.Range("B$firstRow`:B$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Column C - literal ";" terminator, all rows.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = ";"
}

# 4) Column D - concatenation formula, written in two batches (rows
#    185-193 and 194-206) so it mirrors the original workbook's split of
#    the shared-formula range.
$ws.Range("D$firstRow`:D193").Formula = '=A185&" "&B185&" "&C185'
$ws.Range("D194:D$lastRow").Formula = '=A194&" "&B194&" "&C194'

# 5) Update the sheet's selection to the freshly written block, matching
#    the saved view state.
$ws.Range("D$firstRow`:D$lastRow").Select()

Write-Host "Added rows $firstRow-$lastRow for mariana.brider"
